# Estruturação do projeto completa...
#
# The two SharePoint/"Document Information Panel" custom XML parts that were
# shipped with this deck had their package contents swapped:
#   - the part that used to hold the <FormTemplates> (mso-contentType) blob
#     now holds the <ct:contentTypeSchema> blob, and vice-versa
#   - each part's companion item-properties datastore (the one that lists the
#     part's schema references) swaps along with it, so the pairing between a
#     custom XML part and its own schemaRefs stays self-consistent.
# Re-create that by swapping the .XML payloads of the two CustomXMLParts in
# place (their identity / schema references travel with the XML, so this is
# equivalent to relocating the parts).

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$formsNs      = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
$contentTypeNs = "http://schemas.microsoft.com/office/2006/metadata/contentType"

$formsPart = $null
$schemaPart = $null

# Prefer the documented lookup (by root-element namespace); fall back to a
# manual scan of the collection in case SelectByNamespace can't be used.
try {
    $formsScoped = $parts.SelectByNamespace($formsNs)
    if ($formsScoped.Count -ge 1) { $formsPart = $formsScoped.Item(1) }
} catch { }

try {
    $schemaScoped = $parts.SelectByNamespace($contentTypeNs)
    if ($schemaScoped.Count -ge 1) { $schemaPart = $schemaScoped.Item(1) }
} catch { }

if (($formsPart -eq $null) -or ($schemaPart -eq $null)) {
    for ($i = 1; $i -le $parts.Count; $i++) {
        $candidate = $parts.Item($i)
        $xml = $candidate.XML
        if (($formsPart -eq $null) -and ($xml -like "*FormTemplates*")) {
            $formsPart = $candidate
        } elseif (($schemaPart -eq $null) -and ($xml -like "*contentTypeSchema*")) {
            $schemaPart = $candidate
        }
    }
}

# Last-resort fallback: this deck ships the forms blob as customXml/item2.xml
# and the content-type schema blob as customXml/item3.xml, which line up with
# CustomXMLParts.Item(2) / .Item(3) in package order.
if (($formsPart -eq $null) -or ($schemaPart -eq $null)) {
    try {
        if (($formsPart -eq $null) -and ($parts.Count -ge 2)) { $formsPart = $parts.Item(2) }
        if (($schemaPart -eq $null) -and ($parts.Count -ge 3)) { $schemaPart = $parts.Item(3) }
    } catch { }
}

if (($formsPart -ne $null) -and ($schemaPart -ne $null)) {
    $formsXml  = $formsPart.XML
    $schemaXml = $schemaPart.XML

    $formsPart.XML  = $schemaXml
    $schemaPart.XML = $formsXml
}
